# "shiftcal mutants filtering init"
# Adds a new worksheet "ShiftCal - BWD - TWD" after "ShiftCal", fills it
# with the mutants-filtering table, formats it like the existing
# "ShiftCal" sheet, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Reset the current active sheet's selection back to A1 (mirrors the
#    diff: ShiftCal's tabSelected goes false / selection resets to A1
#    once focus moves to the newly-created sheet).
# ---------------------------------------------------------------------
$shiftCal = $wb.Worksheets.Item("ShiftCal")
$shiftCal.Activate()
$shiftCal.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Create the new sheet right after "ShiftCal" and name it.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $shiftCal)
$newSheet.Name = "ShiftCal - BWD - TWD"

# ---------------------------------------------------------------------
# 3. Table data: header + 18 mutant rows.
# ---------------------------------------------------------------------
$fab = "com.google.android.material.floatingactionbutton.FloatingActionButton"

$data = @(
    @("Mutant", "Resource", "Tag", "Operator", "Valid"),
    @("mutant_1", "activity_shift_creator.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_2", "content_shift_creator.xml", "EditText widget deletion", "EditText", "TBD"),
    @("mutant_3", "content_shift_creator.xml", "EditText widget deletion", "EditText", "TBD"),
    @("mutant_4", "content_shift_creator.xml", "Button widget deletion", "android.widget.Button", "TBD"),
    @("mutant_5", "content_shift_creator.xml", "Button widget deletion", "android.widget.Button", "TBD"),
    @("mutant_6", "content_shift_creator.xml", "Button widget deletion", "android.widget.Button", "TBD"),
    @("mutant_7", "activity_calendar.xml", "Button widget deletion", "ImageButton", "TBD"),
    @("mutant_8", "activity_employer_creator.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_9", "content_calendar.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_10", "content_calendar.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_11", "content_about.xml", "Button widget deletion", "Button", "TBD"),
    @("mutant_12", "content_alarm.xml", "EditText widget deletion", "EditText", "TBD"),
    @("mutant_13", "content_alarm.xml", "Button widget deletion", "android.widget.Button", "TBD"),
    @("mutant_14", "content_executed_alarm.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_15", "activity_shifts.xml", "Button widget deletion", $fab, "TBD"),
    @("mutant_16", "content_employer_creator.xml", "EditText widget deletion", "EditText", "TBD"),
    @("mutant_17", "content_theme.xml", "Button widget deletion", "android.widget.Button", "TBD"),
    @("mutant_18", "activity_employers.xml", "Button widget deletion", $fab, "TBD")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# 4. Formatting.
#    - Header row: reuse ShiftCal's bold/bordered header format.
#    - Body columns A, C, E: centered horizontal alignment.
# ---------------------------------------------------------------------
$shiftCal.Range("A1:D1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)
$shiftCal.Range("D1").Copy()
$newSheet.Range("E1").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

$newSheet.Range("A2:A19").HorizontalAlignment = -4108
$newSheet.Range("C2:C19").HorizontalAlignment = -4108
$newSheet.Range("E2:E19").HorizontalAlignment = -4108

# Row heights: header slightly taller, like ShiftCal's header row.
$newSheet.Rows.Item(1).RowHeight = 13.9

# Column widths (A-D custom; E left at sheet default, as in the source).
$newSheet.Columns.Item(1).ColumnWidth = 12.25
$newSheet.Columns.Item(2).ColumnWidth = 23.4167
$newSheet.Columns.Item(3).ColumnWidth = 24.25
$newSheet.Columns.Item(4).ColumnWidth = 23.4167

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 5. Make the new sheet the active tab (3rd sheet, activeTab index 2).
# ---------------------------------------------------------------------
$newSheet.Activate()
